$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")
$tbl = $ws.ListObjects.Item("Assets")

# Grow the Assets table to include 3 new columns (Location, Asset Photo [image], Asset Barcode)
# and 2 new data rows.
$tbl.Resize($ws.Range("B2:I7"))

# New column headers
$ws.Range("G2").Value = "Location"
$ws.Range("H2").Value = "Asset Photo [image]"
$ws.Range("I2").Value = "Asset Barcode"

# Existing rows gain a Location value
$ws.Range("G3").Value = "40.804000,-74.464460"
$ws.Range("G4").Value = "40.804000,-74.464460"
$ws.Range("G5").Value = "40.804000,-74.464460"

# New row 6 - a full asset record
$ws.Range("B6").Value = 123
$ws.Range("C6").Value = "Monitor"
$ws.Range("D6").Value = "Surface Go"
$ws.Range("E6").Value = "256GB HD, 8GB RAM"
$ws.Range("F6").Value = "320QlCiydlQ"
$ws.Range("G6").Value = "52.4804096,-2.146304"

# New row 7 - a partially-filled record (barcode/location/photo captured, rest pending)
$ws.Range("F7").Value = "e2iK7AowT7A"
$ws.Range("G7").Value = "52.45075178,-1.7336544"
$ws.Range("H7").Value = ".\Assets_images\adfbff03195a43cfa06c295e07952f40.png"
